# Update "g19.10" data table: each region's 10-year rolling window shifts
# forward by one year (2025 added, earliest year dropped), so most rows'
# Ano/Valor (and sometimes Posicao/Faltam-dados) values move down one slot
# within their region block; three brand-new rows are appended at the end
# for Sergipe/2025 (and the two following "shift" rows).
# "Chg" lists exactly which columns (A=Regiao,B=Ano,C=Variavel,D=Valor,
# E=Posicao,F=Faltam) actually change for that row, so we only touch cells
# that truly need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{R=2; A="Brasil"; B="01/01/2015"; C="Furto de veículo"; D=88.79597901387544; E=$null; F=1; Chg="D"}
  @{R=3; A="Brasil"; B="01/01/2016"; C="Furto de veículo"; D=95.93705329620661; E=$null; F=1; Chg="D"}
  @{R=4; A="Brasil"; B="01/01/2017"; C="Furto de veículo"; D=95.60364377088423; E=$null; F=1; Chg="D"}
  @{R=5; A="Brasil"; B="01/01/2018"; C="Furto de veículo"; D=92.59800655189056; E=$null; F=1; Chg="D"}
  @{R=6; A="Brasil"; B="01/01/2019"; C="Furto de veículo"; D=89.79331510582048; E=$null; F=1; Chg="D"}
  @{R=7; A="Brasil"; B="01/01/2020"; C="Furto de veículo"; D=71.11763144273382; E=$null; F=1; Chg="D"}
  @{R=8; A="Brasil"; B="01/01/2021"; C="Furto de veículo"; D=72.6315573823851; E=$null; F=0; Chg="D"}
  @{R=9; A="Brasil"; B="01/01/2022"; C="Furto de veículo"; D=82.1331061684723; E=$null; F=0; Chg="D"}
  @{R=10; A="Brasil"; B="01/01/2023"; C="Furto de veículo"; D=75.98479347414431; E=$null; F=0; Chg="D"}
  @{R=11; A="Brasil"; B="01/01/2024"; C="Furto de veículo"; D=71.75411851343904; E=$null; F=0; Chg="D"}
  @{R=12; A="Brasil"; B="01/01/2025"; C="Furto de veículo"; D=39.10563104945712; E=$null; F=0; Chg="ABDF"}
  @{R=13; A="Nordeste"; B="01/01/2015"; C="Furto de veículo"; D=37.21348807397115; E=$null; F=1; Chg="BD"}
  @{R=14; A="Nordeste"; B="01/01/2016"; C="Furto de veículo"; D=47.48042089535613; E=$null; F=1; Chg="BD"}
  @{R=15; A="Nordeste"; B="01/01/2017"; C="Furto de veículo"; D=48.0777904223895; E=$null; F=1; Chg="BD"}
  @{R=16; A="Nordeste"; B="01/01/2018"; C="Furto de veículo"; D=45.71378107315596; E=$null; F=1; Chg="BD"}
  @{R=17; A="Nordeste"; B="01/01/2019"; C="Furto de veículo"; D=45.24399645117818; E=$null; F=1; Chg="BD"}
  @{R=18; A="Nordeste"; B="01/01/2020"; C="Furto de veículo"; D=39.74359243075034; E=$null; F=1; Chg="BDF"}
  @{R=19; A="Nordeste"; B="01/01/2021"; C="Furto de veículo"; D=40.00361392350651; E=$null; F=0; Chg="BD"}
  @{R=20; A="Nordeste"; B="01/01/2022"; C="Furto de veículo"; D=53.76058777621977; E=$null; F=0; Chg="BD"}
  @{R=21; A="Nordeste"; B="01/01/2023"; C="Furto de veículo"; D=55.14400851335632; E=$null; F=0; Chg="BD"}
  @{R=22; A="Nordeste"; B="01/01/2024"; C="Furto de veículo"; D=51.39208374376279; E=$null; F=0; Chg="ABDEF"}
  @{R=23; A="Nordeste"; B="01/01/2025"; C="Furto de veículo"; D=29.92806310753727; E=$null; F=0; Chg="ABDEF"}
  @{R=24; A="Sergipe"; B="01/01/2015"; C="Furto de veículo"; D=36.2708348652949; E=18; F=1; Chg="BDE"}
  @{R=25; A="Sergipe"; B="01/01/2016"; C="Furto de veículo"; D=50.10880001716483; E=18; F=1; Chg="BDE"}
  @{R=26; A="Sergipe"; B="01/01/2017"; C="Furto de veículo"; D=38.01058448583375; E=22; F=1; Chg="BD"}
  @{R=27; A="Sergipe"; B="01/01/2018"; C="Furto de veículo"; D=29.05665081279616; E=25; F=1; Chg="BDE"}
  @{R=28; A="Sergipe"; B="01/01/2019"; C="Furto de veículo"; D=37.45266218394695; E=22; F=1; Chg="BDEF"}
  @{R=29; A="Sergipe"; B="01/01/2020"; C="Furto de veículo"; D=36.26513131341008; E=21; F=1; Chg="BDEF"}
  @{R=30; A="Sergipe"; B="01/01/2021"; C="Furto de veículo"; D=28.47750533632532; E=26; F=0; Chg="BDE"}
  @{R=31; A="Sergipe"; B="01/01/2022"; C="Furto de veículo"; D=38.50982254945205; E=26; F=0; Chg="BD"}
  @{R=32; A="Sergipe"; B="01/01/2023"; C="Furto de veículo"; D=36.06223913262109; E=26; F=0; Chg="ABCDEF"}
  @{R=33; A="Sergipe"; B="01/01/2024"; C="Furto de veículo"; D=32.19900404429514; E=26; F=0; Chg="ABCDEF"}
  @{R=34; A="Sergipe"; B="01/01/2025"; C="Furto de veículo"; D=22.30586229571271; E=25; F=0; Chg="ABCDEF"}
)

foreach ($row in $rows) {
  $r = $row.R
  $chg = $row.Chg

  if ($chg.Contains("A")) {
    $ws.Cells.Item($r, 1).Value = $row.A
  }

  if ($chg.Contains("B")) {
    # Dates are stored as plain text (e.g. "01/01/2025") in this sheet, not
    # as real Excel dates. Force text formatting first so the COM layer
    # doesn't auto-convert the string into a date serial number, then drop
    # back to the default "Normal" style so no stray number format sticks.
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $row.B
    $cell.Style = "Normal"
  }

  if ($chg.Contains("C")) {
    $ws.Cells.Item($r, 3).Value = $row.C
  }

  if ($chg.Contains("D")) {
    $ws.Cells.Item($r, 4).Value = $row.D
  }

  if ($chg.Contains("E")) {
    if ($null -eq $row.E) {
      $ws.Cells.Item($r, 5).ClearContents()
    } else {
      $ws.Cells.Item($r, 5).Value = $row.E
    }
  }

  if ($chg.Contains("F")) {
    $ws.Cells.Item($r, 6).Value = [bool]$row.F
  }
}
